# Trade #6 closed at 2026-02-17 07:57:55 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.95
$summary.Range("B4").Value = -0.05
$summary.Range("B5").Value = -0.17
$summary.Range("B6").Value = 6
$summary.Range("B7").Value = 2
$summary.Range("B9").Value = 33.33

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.95
$status.Range("D4").Value = 6
$status.Range("E4").Value = -0.05
$status.Range("F4").Value = -0.05
$status.Range("G4").Value = 33.33

# --- New trade row data, shared by "All Trades" and "MarketMaking" sheets ---
function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = 6
    # Force the date-looking text to stay a plain string (not auto-converted
    # to a date serial) the same way the existing rows store it, then strip
    # the temporary text-format style so no stray cell style sticks around.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).ClearFormats()
    $ws.Cells.Item($row, 3).Value = "07:57:48"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.79
    $ws.Cells.Item($row, 7).Value = 0.8
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 1.2658
    $ws.Cells.Item($row, 10).Value = 0.01
    $ws.Cells.Item($row, 11).Value = 99.95
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 7

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 7
